{"js": "// The document's Title / Author / Abstract paragraphs were each split into\n// many small runs (one run per word, plus separate runs for the spaces\n// between them). This edit collapses each of those paragraphs back down to\n// a single run carrying the full paragraph text (the visible text itself is\n// unchanged).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\n// Group paragraphs by their style's display name so we find the right\n// paragraph regardless of its position in the document.\nconst byStyle = {};\nfor (const p of paragraphs.items) {\n  if (!byStyle[p.style]) {\n    byStyle[p.style] = [];\n  }\n  byStyle[p.style].push(p);\n}\n\nconst titleText = \"Questions: Arithmetic on complex numbers\";\nconst authorText = \"Charlotte McCarthy\";\nconst abstractText =\n  \"A selection of questions for the study guide on arithmetic on complex numbers.\";\n\nconst titlePara = byStyle[\"Title\"][0];\nconst authorPara = byStyle[\"Author\"][0];\nconst abstractPara = byStyle[\"Abstract\"][0];\n\n// Replacing the whole paragraph range's text merges the (previously\n// word-by-word) runs into a single run.\ntitlePara.insertText(titleText, Word.InsertLocation.replace);\nauthorPara.insertText(authorText, Word.InsertLocation.replace);\nabstractPara.insertText(abstractText, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The document's Title / Author / Abstract paragraphs were each split into\n# many small runs (one run per word, plus separate runs for the spaces\n# between them). This edit collapses each of those paragraphs back down to\n# a single run carrying the full paragraph text (the visible text itself is\n# unchanged).\n\n$d = $word.ActiveDocument\n\n# Map style display name -> the (unchanged) full text that paragraph should\n# end up containing in a single run.\n$targets = @{\n    \"Title\"    = \"Questions: Arithmetic on complex numbers\";\n    \"Author\"   = \"Charlotte McCarthy\";\n    \"Abstract\" = \"A selection of questions for the study guide on arithmetic on complex numbers.\";\n}\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Style.NameLocal\n    if ($targets.ContainsKey($styleName)) {\n        $finalText = $targets[$styleName]\n\n        # Range of the paragraph's text, excluding the trailing paragraph mark.\n        $r = $p.Range\n        $r.MoveEnd(1, -1) | Out-Null\n\n        # Word COM no-ops a Range.Text assignment whose value already equals\n        # the range's (multi-run) text, so first stamp a placeholder to force\n        # the run rewrite, then set the real text - this collapses every run\n        # in the paragraph into a single run.\n        $r.Text = \"~\"\n        $r2 = $p.Range\n        $r2.MoveEnd(1, -1) | Out-Null\n        $r2.Text = $finalText\n    }\n}\n"}
